$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-14 Thursday" "2024-03-15 Friday"

Replace-Text "202÷4=50, 2" "780÷3=260, 0"
Replace-Text "509÷2=254, 1" "910÷2=455, 0"
Replace-Text "695÷7=99, 2" "958÷6=159, 4"
Replace-Text "269÷5=53, 4" "962÷4=240, 2"
Replace-Text "130÷2=65, 0" "205÷2=102, 1"
Replace-Text "297÷8=37, 1" "165÷6=27, 3"
Replace-Text "721÷2=360, 1" "229÷4=57, 1"
Replace-Text "582÷4=145, 2" "755÷6=125, 5"
Replace-Text "294÷4=73, 2" "965÷2=482, 1"
Replace-Text "218÷7=31, 1" "340÷2=170, 0"
Replace-Text "580÷2=290, 0" "175÷5=35, 0"
Replace-Text "678÷6=113, 0" "144÷9=16, 0"
Replace-Text "341÷9=37, 8" "892÷5=178, 2"
Replace-Text "936÷5=187, 1" "864÷8=108, 0"
Replace-Text "572÷6=95, 2" "281÷5=56, 1"
Replace-Text "734÷2=367, 0" "234÷3=78, 0"
Replace-Text "993÷2=496, 1" "453÷2=226, 1"
Replace-Text "523÷7=74, 5" "610÷6=101, 4"
Replace-Text "723÷9=80, 3" "720÷9=80, 0"
Replace-Text "789÷5=157, 4" "833÷7=119, 0"
Replace-Text "828÷2=414, 0" "709÷7=101, 2"
Replace-Text "487÷9=54, 1" "656÷2=328, 0"
Replace-Text "989÷7=141, 2" "101÷8=12, 5"
Replace-Text "975÷9=108, 3" "205÷3=68, 1"
Replace-Text "312÷9=34, 6" "166÷6=27, 4"
